# Update the "Förändrad" (Changed) date column (column C) for every data
# row (rows 2-171) from 45189 (2023-09-20) to 45190 (2023-09-21).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 171; $r++) {
    $ws.Cells.Item($r, 3).Value = 45190
}
